# "Generate Report for Handback"
#
# This script updates the localization-status workbook to reflect that the
# handback (target -> en-US sync) step has completed for both the zh-cn and
# de-de languages:
#   - the "Status" column switches from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it appears
#   - the zh-cn and de-de report sheets get their "Latest Target File",
#     "Latest Handback File" and "Latest Handback DateTime" columns filled
#     in (including real hyperlinks for the target file names)
#   - a handful of columns are widened so the new, longer values are legible

$wb = $excel.ActiveWorkbook

$Overview = $wb.Worksheets.Item(1)   # "Overview"
$ZhCn     = $wb.Worksheets.Item(2)   # "zh-cn"
$DeDe     = $wb.Worksheets.Item(3)   # "de-de"

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Update every cell currently showing the old status text so the
#    shared string is fully replaced by the new status text.
# ---------------------------------------------------------------------
$Overview.Range("E2").Value = $statusNew
$Overview.Range("F2").Value = $statusNew
$Overview.Range("E3").Value = $statusNew
$Overview.Range("F3").Value = $statusNew

$ZhCn.Range("C2").Value = $statusNew
$ZhCn.Range("C3").Value = $statusNew

$DeDe.Range("C2").Value = $statusNew
$DeDe.Range("C3").Value = $statusNew

# ---------------------------------------------------------------------
# 2. zh-cn sheet: fill in Latest Target File (I), Latest Handback File (J)
#    and Latest Handback DateTime (K) for both data rows.
# ---------------------------------------------------------------------
$ZhCn.Hyperlinks.Add(
    $ZhCn.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/oltest/blob/5ab631ea6f0732e30e27ee0ec371ea7e83c88c6d/e2e/60ddc312-19a7-4e79-8f67-60806d417789.md",
    "",
    "",
    "60ddc312-19a7-4e79-8f67-60806d417789.md"
) | Out-Null
$ZhCn.Range("J2").Value = "60ddc312-19a7-4e79-8f67-60806d417789.2dee38ea3909e30db3d71c188faabfd05e858e71.zh-cn.xlf"
$ZhCn.Range("K2").Value = "2016-08-13 04:58:12"

$ZhCn.Hyperlinks.Add(
    $ZhCn.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/oltest/blob/5ab631ea6f0732e30e27ee0ec371ea7e83c88c6d/e2e/dbca198b-97a7-4d9e-9af7-45a3899e2554.md",
    "",
    "",
    "dbca198b-97a7-4d9e-9af7-45a3899e2554.md"
) | Out-Null
$ZhCn.Range("J3").Value = "dbca198b-97a7-4d9e-9af7-45a3899e2554.1ccd2c5ee60df00c46018291035e1210f3d31446.zh-cn.xlf"
$ZhCn.Range("K3").Value = "2016-08-13 04:58:12"

# ---------------------------------------------------------------------
# 3. de-de sheet: fill in Latest Target File (I), Latest Handback File (J)
#    and Latest Handback DateTime (K) for both data rows.
# ---------------------------------------------------------------------
$DeDe.Hyperlinks.Add(
    $DeDe.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/oltest/blob/5ab631ea6f0732e30e27ee0ec371ea7e83c88c6d/e2e/60ddc312-19a7-4e79-8f67-60806d417789.md",
    "",
    "",
    "60ddc312-19a7-4e79-8f67-60806d417789.md"
) | Out-Null
$DeDe.Range("J2").Value = "60ddc312-19a7-4e79-8f67-60806d417789.2dee38ea3909e30db3d71c188faabfd05e858e71.de-de.xlf"
$DeDe.Range("K2").Value = "2016-08-13 04:58:22"

$DeDe.Hyperlinks.Add(
    $DeDe.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/oltest/blob/5ab631ea6f0732e30e27ee0ec371ea7e83c88c6d/e2e/dbca198b-97a7-4d9e-9af7-45a3899e2554.md",
    "",
    "",
    "dbca198b-97a7-4d9e-9af7-45a3899e2554.md"
) | Out-Null
$DeDe.Range("J3").Value = "dbca198b-97a7-4d9e-9af7-45a3899e2554.1ccd2c5ee60df00c46018291035e1210f3d31446.de-de.xlf"
$DeDe.Range("K3").Value = "2016-08-13 04:58:22"

# ---------------------------------------------------------------------
# 4. Widen columns so the newly-populated / longer text is fully visible.
# ---------------------------------------------------------------------
$Overview.Columns.Item(5).ColumnWidth = 29.15   # E
$Overview.Columns.Item(6).ColumnWidth = 29.15   # F

$ZhCn.Columns.Item(3).ColumnWidth = 29.15        # C - Status
$ZhCn.Columns.Item(9).ColumnWidth = 39.2         # I - Latest Target File
$ZhCn.Columns.Item(10).ColumnWidth = 39.2        # J - Latest Handback File

$DeDe.Columns.Item(3).ColumnWidth = 29.15        # C - Status
$DeDe.Columns.Item(9).ColumnWidth = 39.2         # I - Latest Target File
$DeDe.Columns.Item(10).ColumnWidth = 39.2        # J - Latest Handback File
